# Rename the inline picture shapes that live in the document's headers/
# footers. Two copies of the Pearson Edexcel logo (in the default and
# first-page footers) go from "image1.png" to "image2.png", and the BTEC
# logo (in the first-page header) goes from "image2.jpg" to "image1.jpg".
#
# The pictures aren't reachable via $d.InlineShapes (that collection only
# covers the main body story) - they have to be found through each
# Section's Headers/Footers collections instead.

$d = $word.ActiveDocument

function Rename-LogoShapesIn($range) {
    if ($range.InlineShapes.Count -eq 0) {
        return
    }
    foreach ($shape in $range.InlineShapes) {
        $alt = $shape.AlternativeText
        if ($alt -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
            $shape.Name = "image2.png"
        } elseif ($alt -eq "BTec_Logo-Orange") {
            $shape.Name = "image1.jpg"
        }
    }
}

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le 3; $i++) {
        $hdr = $sec.Headers.Item($i)
        if ($hdr.Exists) {
            Rename-LogoShapesIn $hdr.Range
        }
        $ftr = $sec.Footers.Item($i)
        if ($ftr.Exists) {
            Rename-LogoShapesIn $ftr.Range
        }
    }
}
